$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/4/2025  Through  8/10/2025"

# --- Weekly crime-data table updates (rows 15-31) ---
# Row 15
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 3
$ws.Range("D15").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("G15").Value = "'0"
$ws.Range("D15").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("H15").Value = "'***.*"
$ws.Range("I15").Value = 20
$ws.Range("K15").Value = 42.857142857142
$ws.Range("L15").Value = 66.666666666666
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = -75

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 36.363636363636
$ws.Range("I16").Value = 118
$ws.Range("J16").Value = 121
$ws.Range("K16").Value = -2.479338842975
$ws.Range("L16").Value = 20.408163265306
$ws.Range("M16").Value = -47.555555555555
$ws.Range("N16").Value = -90.90909090909

# Row 17
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1100
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 117.647058823529
$ws.Range("I17").Value = 293
$ws.Range("J17").Value = 206
$ws.Range("K17").Value = 42.233009708737
$ws.Range("L17").Value = 44.334975369458
$ws.Range("M17").Value = 28.508771929824
$ws.Range("N17").Value = -46.823956442831

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 150
$ws.Range("I18").Value = 85
$ws.Range("J18").Value = 72
$ws.Range("K18").Value = 18.055555555555
$ws.Range("L18").Value = 8.974358974358
$ws.Range("M18").Value = -49.704142011834
$ws.Range("N18").Value = -95.002939447383

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 28.571428571428
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -2.857142857142
$ws.Range("I19").Value = 293
$ws.Range("J19").Value = 327
$ws.Range("K19").Value = -10.397553516819
$ws.Range("L19").Value = -11.212121212121
$ws.Range("M19").Value = -22.691292875989
$ws.Range("N19").Value = -55.335365853658

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = -28.571428571428
$ws.Range("I20").Value = 87
$ws.Range("J20").Value = 98
$ws.Range("K20").Value = -11.224489795918
$ws.Range("L20").Value = 38.095238095238
$ws.Range("M20").Value = -26.890756302521
$ws.Range("N20").Value = -94.390715667311

# Row 21
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 112.5
$ws.Range("F21").Value = 114
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = 29.545454545454
$ws.Range("I21").Value = 896
$ws.Range("J21").Value = 840
$ws.Range("K21").Value = 6.666666666666
$ws.Range("L21").Value = 13.705583756345
$ws.Range("M21").Value = -21.5411558669
$ws.Range("N21").Value = -84.720327421555

# Row 22
$ws.Range("F22").Value = 1
$ws.Range("M22").Value = -64.705882352941

# Row 24
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = -20.51282051282
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 135
$ws.Range("H24").Value = -10.37037037037
$ws.Range("I24").Value = 903
$ws.Range("J24").Value = 1061
$ws.Range("K24").Value = -14.891611687087
$ws.Range("L24").Value = -20.017714791851
$ws.Range("M24").Value = 15.91784338896

# Row 25
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 30
$ws.Range("E25").Value = -43.333333333333
$ws.Range("F25").Value = 67
$ws.Range("G25").Value = 95
$ws.Range("H25").Value = -29.473684210526
$ws.Range("I25").Value = 451
$ws.Range("J25").Value = 667
$ws.Range("K25").Value = -32.383808095952
$ws.Range("L25").Value = -31.562974203338

# Row 26
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -21.428571428571
$ws.Range("F26").Value = 51
$ws.Range("G26").Value = 60
$ws.Range("H26").Value = -15
$ws.Range("I26").Value = 408
$ws.Range("J26").Value = 415
$ws.Range("K26").Value = -1.686746987951
$ws.Range("L26").Value = 12.087912087912
$ws.Range("M26").Value = -13.006396588486

# Row 27
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 3
$ws.Range("D27").Copy() | Out-Null
$ws.Range("G27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("G27").Value = "'0"
$ws.Range("D27").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("H27").Value = "'***.*"
$ws.Range("I27").Value = 23
$ws.Range("K27").Value = 15
$ws.Range("L27").Value = 9.523809523809

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 28.571428571428
$ws.Range("I28").Value = 29
$ws.Range("J28").Value = 46
$ws.Range("K28").Value = -36.95652173913
$ws.Range("L28").Value = -29.268292682926

# Row 29
$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 9
$ws.Range("K29").Value = -55.555555555555
$ws.Range("N29").Value = -95.604395604395

# Row 30
$ws.Range("J30").Value = 8
$ws.Range("K30").Value = -50
$ws.Range("N30").Value = -94.444444444444

# Row 31
$ws.Range("C31").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D31").Value = "'0"
$ws.Range("C31").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E31").Value = "'***.*"
$ws.Range("L31").Value = 9.090909090909

